$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell $ws.Cells.Item(2, 4) '29.455.51'
Set-TextCell $ws.Cells.Item(2, 5) '  +0.94%  '
Set-TextCell $ws.Cells.Item(3, 4) '1.918.23'
Set-TextCell $ws.Cells.Item(3, 5) '  +1.64%  '
Set-TextCell $ws.Cells.Item(4, 4) '1.009'
Set-TextCell $ws.Cells.Item(4, 5) '  +0.72%  '
Set-TextCell $ws.Cells.Item(5, 4) '325.46'
Set-TextCell $ws.Cells.Item(5, 5) '  +0.93%  '
Set-TextCell $ws.Cells.Item(6, 4) '1.007'
Set-TextCell $ws.Cells.Item(6, 5) '  +0.55%  '
Set-TextCell $ws.Cells.Item(7, 4) '0.4827'
Set-TextCell $ws.Cells.Item(7, 5) '  +2.76%  '
Set-TextCell $ws.Cells.Item(8, 4) '0.4078'
Set-TextCell $ws.Cells.Item(8, 5) '  +1.34%  '
Set-TextCell $ws.Cells.Item(9, 4) '0.08184'
Set-TextCell $ws.Cells.Item(9, 5) '  +2.24%  '
Set-TextCell $ws.Cells.Item(10, 4) '1.020'
Set-TextCell $ws.Cells.Item(10, 5) '  +2.83%  '
Set-TextCell $ws.Cells.Item(11, 4) '23.47'
Set-TextCell $ws.Cells.Item(11, 5) '  +3.76%  '
Set-TextCell $ws.Cells.Item(12, 4) '1.942.27'
Set-TextCell $ws.Cells.Item(12, 5) '  +2.55%  '
Set-TextCell $ws.Cells.Item(13, 4) '6.036'
Set-TextCell $ws.Cells.Item(13, 5) '  +2.09%  '
Set-TextCell $ws.Cells.Item(14, 4) '7.223'
Set-TextCell $ws.Cells.Item(14, 5) '  +3.10%  '
Set-TextCell $ws.Cells.Item(15, 4) '91.23'
Set-TextCell $ws.Cells.Item(15, 5) '  +2.15%  '
Set-TextCell $ws.Cells.Item(16, 4) '0.06793'
Set-TextCell $ws.Cells.Item(16, 5) '  +2.48%  '
Set-TextCell $ws.Cells.Item(17, 5) '  +0.61%  '
Set-TextCell $ws.Cells.Item(18, 4) '0.00001038'
Set-TextCell $ws.Cells.Item(18, 5) '  +1.40%  '
Set-TextCell $ws.Cells.Item(19, 5) '  +2.28%  '
Set-TextCell $ws.Cells.Item(20, 4) '1.007'
Set-TextCell $ws.Cells.Item(20, 5) '  +0.63%  '
Set-TextCell $ws.Cells.Item(21, 4) '29.493.88'
Set-TextCell $ws.Cells.Item(21, 5) '  +1.08%  '
Set-TextCell $ws.Cells.Item(22, 4) '5.634'
Set-TextCell $ws.Cells.Item(22, 5) '  +2.68%  '
Set-TextCell $ws.Cells.Item(23, 4) '11.78'
Set-TextCell $ws.Cells.Item(23, 5) '  +0.92%  '
Set-TextCell $ws.Cells.Item(24, 4) '2.192'
Set-TextCell $ws.Cells.Item(24, 5) '  +0.75%  '
Set-TextCell $ws.Cells.Item(25, 4) '2.143.88'
Set-TextCell $ws.Cells.Item(25, 5) '  +1.18%  '
Set-TextCell $ws.Cells.Item(26, 4) '6.654'
Set-TextCell $ws.Cells.Item(26, 5) '  +10.92%  '
Set-TextCell $ws.Cells.Item(27, 4) '156.64'
Set-TextCell $ws.Cells.Item(27, 5) '  +1.17%  '
Set-TextCell $ws.Cells.Item(28, 4) '20.05'
Set-TextCell $ws.Cells.Item(28, 5) '  +2.30%  '
Set-TextCell $ws.Cells.Item(29, 4) '2.112'
Set-TextCell $ws.Cells.Item(29, 5) '  +1.63%  '
Set-TextCell $ws.Cells.Item(30, 4) '120.32'
Set-TextCell $ws.Cells.Item(30, 5) '  +2.81%  '
Set-TextCell $ws.Cells.Item(31, 4) '1.021'
Set-TextCell $ws.Cells.Item(31, 5) '  -0.19%  '
Set-TextCell $ws.Cells.Item(32, 4) '0.09564'
Set-TextCell $ws.Cells.Item(32, 5) '  +1.74%  '
Set-TextCell $ws.Cells.Item(33, 4) '5.518'
Set-TextCell $ws.Cells.Item(33, 5) '  +3.20%  '
Set-TextCell $ws.Cells.Item(34, 4) '3.562'
Set-TextCell $ws.Cells.Item(34, 5) '  +0.63%  '
Set-TextCell $ws.Cells.Item(35, 4) '1.382'
Set-TextCell $ws.Cells.Item(35, 5) '  +0.50%  '
Set-TextCell $ws.Cells.Item(36, 4) '0.02284'
Set-TextCell $ws.Cells.Item(36, 5) '  +2.19%  '
Set-TextCell $ws.Cells.Item(37, 4) '0.06129'
Set-TextCell $ws.Cells.Item(37, 5) '  +1.54%  '
Set-TextCell $ws.Cells.Item(38, 4) '1.180'
Set-TextCell $ws.Cells.Item(38, 5) '  +0.89%  '
Set-TextCell $ws.Cells.Item(39, 4) '0.5978'
Set-TextCell $ws.Cells.Item(39, 5) '  +2.97%  '
Set-TextCell $ws.Cells.Item(40, 4) '8.031'
Set-TextCell $ws.Cells.Item(40, 5) '  +0.57%  '
Set-TextCell $ws.Cells.Item(41, 4) '10.80'
Set-TextCell $ws.Cells.Item(41, 5) '  +7.90%  '
Set-TextCell $ws.Cells.Item(42, 4) '0.1856'
Set-TextCell $ws.Cells.Item(42, 5) '  +1.66%  '
Set-TextCell $ws.Cells.Item(43, 4) '1.282'
Set-TextCell $ws.Cells.Item(43, 5) '  +0.80%  '
Set-TextCell $ws.Cells.Item(44, 5) '  -1.14%  '
Set-TextCell $ws.Cells.Item(45, 4) '0.07606'
Set-TextCell $ws.Cells.Item(46, 4) '12.42'
Set-TextCell $ws.Cells.Item(46, 5) '  +2.36%  '
Set-TextCell $ws.Cells.Item(47, 4) '0.5570'
Set-TextCell $ws.Cells.Item(47, 5) '  +2.13%  '
Set-TextCell $ws.Cells.Item(48, 4) '1.956'
Set-TextCell $ws.Cells.Item(48, 5) '  +3.06%  '
Set-TextCell $ws.Cells.Item(49, 4) '117.41'
Set-TextCell $ws.Cells.Item(49, 5) '  +3.68%  '
Set-TextCell $ws.Cells.Item(50, 5) '  +4.56%  '
Set-TextCell $ws.Cells.Item(51, 4) '72.67'
Set-TextCell $ws.Cells.Item(51, 5) '  +2.41%  '
